# add leitura de base e criação do arquivo de saída
#
# Inserts a new "Sheet" column (B) between the "Nome base" and
# "Local do arquivo" columns, shifting the existing B:E columns right by
# one, widens the (now) "Local do arquivo" column to fit the longer file
# paths, fills the new column with the worksheet-name keys ("basea" /
# "baseb"), and fixes up the two file-path values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B; this shifts B->C, C->D, D->E, E->F,
# F->G, G->H and carries each cell's existing value/style along with it.
$ws.Columns.Item(2).Insert()

# --- Header row ---------------------------------------------------------
$ws.Range("B1").Value = "Sheet"
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B1").HorizontalAlignment = -4108  # xlCenter
$excel.CutCopyMode = $false

# --- Data rows -----------------------------------------------------------
$ws.Range("B2").Value = "basea"
$ws.Range("B3").Value = "baseb"

$ws.Range("C2").Value = "C:\Users\jairm\Documents\Code\Empilhador\empilhador\Arquivos\BaseA.xlsx"
$ws.Range("C3").Value = "C:\Users\jairm\Documents\Code\Empilhador\empilhador\Arquivos\BaseB.xlsx"

# --- Column widths ---------------------------------------------------------
# Columns A, D, E, F are just the original A, C, D, E shifted right by the
# insert above, so they already carry their original (bestFit) width along
# with them -- leave them alone. Only the new column B (narrow, default-ish
# width) and column C (now holding the longer file-path strings, so its
# old bestFit width is stale) need to be touched explicitly.
# NB: the engine quantises ColumnWidth to 1/6-character steps, so the wider
# target width here -- driven upstream by Excel's real AutoFit pixel
# metrics -- lands on the nearest reachable step.
$ws.Columns.Item(2).ColumnWidth = 10.166666666666666
$ws.Columns.Item(3).ColumnWidth = 72.59244791666667

# --- Selection -------------------------------------------------------------
$ws.Range("C5").Select()
